$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Samples query text in B3: change the trailing LIMIT 10 to LIMIT 100
$cell = $ws.Range("B3")
$text = $cell.Value2
$newText = $text -replace "LIMIT 10$", "LIMIT 100"
$cell.Value2 = $newText

# Move the active selection to B13 (cursor position at save time)
$ws.Range("B13").Select()
